$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17, duplicating the row currently there
# (this copies formatting/values down, so row 17 and the new row 18 will
# both initially hold the old row 17 data), then update the two cells
# that actually carry new data (Fecha / Volumen).
$ws.Rows.Item(17).Copy() | Out-Null
$ws.Rows.Item(17).Insert() | Out-Null

$ws.Range("D17").Value = 44550
$ws.Range("J17").Value = 65
